$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (11th column), shifting existing
# columns K.. to the right.
$ws.Columns.Item(11).Insert()

# Set header for the newly inserted column K
$ws.Cells.Item(1, 11).Value = "NPC"

# Update the selection to match the target state
$ws.Range("I12").Select()
